# "result=yes/no, added % symbol in percentage"
#
# The "probability" column (M) on Sheet1 currently stores plain decimal
# fraction text such as "0.92". This rewrites each of those cells as a
# percentage-formatted text string such as "92.00%" (a handful of rows
# carry more precise, non-rounded percentages in the source data, e.g.
# "63.76%", "60.83%", "89.33%", "51.50%", "76.80%").
#
# Data layout below: "<row number>|<new text value for column M>".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @"
2|92.00%
3|91.00%
4|73.00%
5|54.00%
6|71.00%
7|75.00%
8|75.00%
9|67.00%
10|67.00%
11|76.00%
12|91.00%
13|98.00%
14|93.00%
15|92.00%
16|91.00%
17|73.00%
18|81.00%
19|76.00%
20|63.76%
21|60.83%
22|66.00%
23|83.00%
24|96.00%
25|75.00%
26|61.00%
27|54.00%
28|69.00%
29|84.00%
30|56.00%
31|76.00%
32|86.00%
33|86.00%
34|82.00%
35|93.00%
36|68.00%
37|63.00%
38|85.00%
39|74.00%
40|89.33%
41|70.00%
42|62.00%
43|90.00%
44|66.00%
45|93.00%
46|61.00%
47|97.00%
48|67.00%
49|51.00%
50|97.00%
51|70.00%
52|85.00%
53|89.00%
54|85.00%
55|61.00%
56|69.00%
57|75.00%
58|93.00%
59|90.00%
60|58.00%
61|86.00%
62|72.00%
63|51.00%
64|99.00%
65|70.00%
66|51.50%
67|87.00%
68|90.00%
69|54.00%
70|53.00%
71|84.00%
72|57.00%
73|96.00%
74|71.00%
75|58.00%
76|90.00%
77|66.00%
78|58.00%
79|84.00%
80|74.00%
81|76.80%
82|76.00%
83|92.00%
84|57.00%
85|69.00%
86|82.00%
87|87.00%
88|88.00%
89|62.00%
90|68.00%
91|99.00%
92|78.00%
93|62.00%
94|90.00%
95|63.00%
96|65.00%
97|62.00%
98|77.00%
99|77.00%
100|82.00%
101|63.00%
102|82.00%
103|94.00%
104|62.00%
105|78.00%
106|94.00%
107|76.00%
108|55.00%
109|66.00%
110|70.00%
111|72.00%
"@

$lines = $data -split "`n"

$targetRange = $ws.Range("M2:M111")

# Force the cells to stay plain text so a value like "92.00%" is stored
# literally instead of being reinterpreted as the number 0.92 with a
# percentage number format.
$targetRange.NumberFormat = "@"

foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|', 2
    $row = [int]$parts[0]
    $text = $parts[1]

    $ws.Cells.Item($row, 13).Value = $text   # column M = 13 ("probability")
}

# Restore the default style so the cells end up identical (no explicit
# style index) to how they looked before this text-format workaround.
$targetRange.Style = "Normal"

Write-Host "Updated $($lines.Count) probability cells in column M."
